# error solve ifrs list
# Update the financial figures (rows 2-9) on the "company_list" sheet so the
# values match the corrected ("error solve") dataset. A handful of cells in
# rows 7-9 (the analyst-estimate years) also drop a trailing, now-unused
# metric column (N/Y/AE/AF) entirely, so those are cleared rather than set.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 1760
$ws.Range("E2").Value = 68
$ws.Range("F2").Value = 68
$ws.Range("G2").Value = 45
$ws.Range("H2").Value = 31
$ws.Range("I2").Value = 20
$ws.Range("J2").Value = 11
$ws.Range("K2").Value = 1416
$ws.Range("L2").Value = 642
$ws.Range("M2").Value = 774
$ws.Range("N2").Value = 742
$ws.Range("O2").Value = 32
$ws.Range("P2").Value = 51
$ws.Range("Q2").Value = 50
$ws.Range("R2").Value = -91
$ws.Range("S2").Value = -3
$ws.Range("T2").Value = 85
$ws.Range("U2").Value = -35
$ws.Range("V2").Value = 292
$ws.Range("W2").Value = 3.84
$ws.Range("X2").Value = 1.79
$ws.Range("Y2").Value = 2.72
$ws.Range("Z2").Value = 2.27
$ws.Range("AA2").Value = 82.91
$ws.Range("AB2").Value = 1620.35
$ws.Range("AC2").Value = 411
$ws.Range("AD2").Value = 36.04
$ws.Range("AE2").Value = 15230
$ws.Range("AF2").Value = 0.97
$ws.Range("AG2").Value = 250
$ws.Range("AH2").Value = 1.69
$ws.Range("AI2").Value = 60.83
$ws.Range("AJ2").Value = 4876849

# Row 3
$ws.Range("D3").Value = 1808
$ws.Range("E3").Value = 213
$ws.Range("F3").Value = 213
$ws.Range("G3").Value = 196
$ws.Range("H3").Value = 118
$ws.Range("I3").Value = 97
$ws.Range("J3").Value = 22
$ws.Range("K3").Value = 1588
$ws.Range("L3").Value = 673
$ws.Range("M3").Value = 915
$ws.Range("N3").Value = 873
$ws.Range("O3").Value = 42
$ws.Range("P3").Value = 53
$ws.Range("Q3").Value = 163
$ws.Range("R3").Value = -125
$ws.Range("S3").Value = 56
$ws.Range("T3").Value = 71
$ws.Range("U3").Value = 93
$ws.Range("V3").Value = 325
$ws.Range("W3").Value = 11.76
$ws.Range("X3").Value = 6.54
$ws.Range("Y3").Value = 11.96
$ws.Range("Z3").Value = 7.88
$ws.Range("AA3").Value = 73.55
$ws.Range("AB3").Value = 1777.12
$ws.Range("AC3").Value = 1951
$ws.Range("AD3").Value = 15.92
$ws.Range("AE3").Value = 16956
$ws.Range("AF3").Value = 1.83
$ws.Range("AG3").Value = 500
$ws.Range("AH3").Value = 1.61
$ws.Range("AI3").Value = 26.65
$ws.Range("AJ3").Value = 5152078

# Row 4
$ws.Range("D4").Value = 1768
$ws.Range("E4").Value = 205
$ws.Range("F4").Value = 205
$ws.Range("G4").Value = 212
$ws.Range("H4").Value = 153
$ws.Range("I4").Value = 135
$ws.Range("J4").Value = 18
$ws.Range("K4").Value = 1972
$ws.Range("L4").Value = 861
$ws.Range("M4").Value = 1111
$ws.Range("N4").Value = 1007
$ws.Range("O4").Value = 104
$ws.Range("P4").Value = 55
$ws.Range("Q4").Value = 262
$ws.Range("R4").Value = -302
$ws.Range("S4").Value = 39
$ws.Range("T4").Value = 33
$ws.Range("U4").Value = 229
$ws.Range("V4").Value = 389
$ws.Range("W4").Value = 11.58
$ws.Range("X4").Value = 8.630000000000001
$ws.Range("Y4").Value = 14.35
$ws.Range("Z4").Value = 8.57
$ws.Range("AA4").Value = 77.47
$ws.Range("AB4").Value = 1981.84
$ws.Range("AC4").Value = 2547
$ws.Range("AD4").Value = 7.6
$ws.Range("AE4").Value = 19178
$ws.Range("AF4").Value = 1.01
$ws.Range("AG4").Value = 500
$ws.Range("AH4").Value = 2.58
$ws.Range("AI4").Value = 19.46
$ws.Range("AJ4").Value = 5304983

# Row 5
$ws.Range("D5").Value = 2043
$ws.Range("E5").Value = 82
$ws.Range("F5").Value = 82
$ws.Range("G5").Value = 57
$ws.Range("H5").Value = 51
$ws.Range("I5").Value = 52
$ws.Range("J5").Value = -1
$ws.Range("K5").Value = 1902
$ws.Range("L5").Value = 798
$ws.Range("M5").Value = 1104
$ws.Range("N5").Value = 1017
$ws.Range("O5").Value = 86
$ws.Range("P5").Value = 55
$ws.Range("Q5").Value = 78
$ws.Range("R5").Value = 61
$ws.Range("S5").Value = -179
$ws.Range("T5").Value = 31
$ws.Range("U5").Value = 47
$ws.Range("V5").Value = 298
$ws.Range("W5").Value = 4.02
$ws.Range("X5").Value = 2.48
$ws.Range("Y5").Value = 5.14
$ws.Range("Z5").Value = 2.62
$ws.Range("AA5").Value = 72.29000000000001
$ws.Range("AB5").Value = 1762.99
$ws.Range("AC5").Value = 979
$ws.Range("AD5").Value = 16.85
$ws.Range("AE5").Value = 19259
$ws.Range("AF5").Value = 0.86
$ws.Range("AG5").Value = 300
$ws.Range("AH5").Value = 1.82
$ws.Range("AI5").Value = 30.47
$ws.Range("AJ5").Value = 5335564

# Row 6
$ws.Range("D6").Value = 1914
$ws.Range("E6").Value = 68
$ws.Range("F6").Value = 68
$ws.Range("G6").Value = 85
$ws.Range("H6").Value = 54
$ws.Range("I6").Value = 58
$ws.Range("K6").Value = 1841
$ws.Range("L6").Value = 712
$ws.Range("M6").Value = 1129
$ws.Range("N6").Value = 1058
$ws.Range("P6").Value = 55
$ws.Range("Q6").Value = 102
$ws.Range("R6").Value = -17
$ws.Range("S6").Value = -70
$ws.Range("T6").Value = 19
$ws.Range("U6").Value = 83
$ws.Range("V6").Value = 254
$ws.Range("W6").Value = 3.54
$ws.Range("X6").Value = 2.8
$ws.Range("Y6").Value = 5.62
$ws.Range("Z6").Value = 2.87
$ws.Range("AA6").Value = 63.05
$ws.Range("AB6").Value = 1831.29
$ws.Range("AC6").Value = 1094
$ws.Range("AD6").Value = 12.57
$ws.Range("AE6").Value = 20023
$ws.Range("AF6").Value = 0.6899999999999999
$ws.Range("AG6").Value = 500
$ws.Range("AH6").Value = 3.64
$ws.Range("AI6").Value = 45.26
$ws.Range("AJ6").Value = 5335564

# Row 7
$ws.Range("D7").Value = 2014
$ws.Range("E7").Value = 262
$ws.Range("G7").Value = 279
$ws.Range("H7").Value = 193
$ws.Range("I7").Value = 164
$ws.Range("K7").Value = 2193
$ws.Range("L7").Value = 782
$ws.Range("M7").Value = 1411
$ws.Range("P7").Value = 55
$ws.Range("Q7").Value = 405
$ws.Range("R7").Value = -58
$ws.Range("S7").Value = -89
$ws.Range("T7").Value = 15
$ws.Range("U7").Value = 432
$ws.Range("W7").Value = 13.01
$ws.Range("X7").Value = 9.58
$ws.Range("Z7").Value = 9.57
$ws.Range("AA7").Value = 55.42
$ws.Range("AC7").Value = 3074
$ws.Range("AD7").Value = 5.51
$ws.Range("AG7").Value = 500
$ws.Range("AH7").Value = 2.95
$ws.Range("AI7").Value = 16.27
$ws.Range("AE7,AF7,N7,Y7").ClearContents()

# Row 8
$ws.Range("D8").Value = 2093
$ws.Range("E8").Value = 274
$ws.Range("G8").Value = 283
$ws.Range("H8").Value = 179
$ws.Range("I8").Value = 156
$ws.Range("K8").Value = 2299
$ws.Range("L8").Value = 759
$ws.Range("M8").Value = 1540
$ws.Range("P8").Value = 55
$ws.Range("Q8").Value = 201
$ws.Range("R8").Value = -36
$ws.Range("S8").Value = -67
$ws.Range("T8").Value = 21
$ws.Range("U8").Value = 180
$ws.Range("W8").Value = 13.09
$ws.Range("X8").Value = 8.550000000000001
$ws.Range("Z8").Value = 7.97
$ws.Range("AA8").Value = 49.29
$ws.Range("AC8").Value = 2924
$ws.Range("AD8").Value = 5.8
$ws.Range("AG8").Value = 500
$ws.Range("AH8").Value = 2.95
$ws.Range("AI8").Value = 17.1
$ws.Range("AE8,AF8,N8,Y8").ClearContents()

# Row 9
$ws.Range("D9").Value = 2114
$ws.Range("E9").Value = 279
$ws.Range("G9").Value = 289
$ws.Range("H9").Value = 183
$ws.Range("I9").Value = 159
$ws.Range("K9").Value = 2402
$ws.Range("L9").Value = 730
$ws.Range("M9").Value = 1672
$ws.Range("P9").Value = 55
$ws.Range("Q9").Value = 206
$ws.Range("R9").Value = -37
$ws.Range("S9").Value = -67
$ws.Range("T9").Value = 21
$ws.Range("U9").Value = 185
$ws.Range("W9").Value = 13.2
$ws.Range("X9").Value = 8.66
$ws.Range("Z9").Value = 7.79
$ws.Range("AA9").Value = 43.66
$ws.Range("AC9").Value = 2980
$ws.Range("AD9").Value = 5.69
$ws.Range("AG9").Value = 500
$ws.Range("AH9").Value = 2.95
$ws.Range("AI9").Value = 16.78
$ws.Range("AE9,AF9,N9,Y9").ClearContents()

